$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShapeImage($range, $newName) {
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $inlineShape = $range.InlineShapes.Item($i)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        $shape.ConvertToInlineShape()
    }
}

# Footer (primary / default) -> PearsonLogo: image1.png -> image2.png
$footerPrimary = $sec.Footers.Item(1)
Rename-InlineShapeImage $footerPrimary.Range "image2.png"

# Footer (first page) -> PearsonLogo: image1.png -> image2.png
$footerFirst = $sec.Footers.Item(2)
Rename-InlineShapeImage $footerFirst.Range "image2.png"

# Header (first page) -> BTec_Logo-Orange: image2.jpg -> image1.jpg
$headerFirst = $sec.Headers.Item(2)
Rename-InlineShapeImage $headerFirst.Range "image1.jpg"
